$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Added goal weight into the velocity-adjustment totals (row 39):
#  - C39 (x): previously SUM(D2:D5)/4, now subtracts -1 (adds the goal's x weight)
#  - D39 (y): previously SUM(E2:E5)/4, now subtracts 2 (adds the goal's y weight)
$ws.Range("C39").Formula = "=SUM(D2:D5)/4-(-1)"
$ws.Range("D39").Formula = "=SUM(E2:E5)/4-2"

# Highlight the goal-adjustment x value (D42) with a 2-decimal number format
$ws.Range("D42").NumberFormat = "0.00"

# Move the active selection to the goal-adjustment row
$ws.Range("C42").Select()
